$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2
$ws.Range("A2").Value = "ECs"
$ws.Range("B2").Value = "Il7"
$ws.Range("C2").Value = "Il7r"
$ws.Range("D2").Value = "ECs"
$ws.Range("E2").Value = 1
$ws.Range("F2").Value = 0.3333333333333333
$ws.Range("G2").Value = 0.1621203333333333
$ws.Range("H2").Value = 0.486361
$ws.Range("I2").Value = 0.1481290732860657
$ws.Range("J2").Value = 0.1532739727236171
$ws.Range("K2").Value = 1
$ws.Range("L2").Value = 0.3333333333333333
$ws.Range("M2").Value = 0.115185
$ws.Range("N2").Value = 0.345555
$ws.Range("O2").Value = 0.001968886112489112
$ws.Range("P2").Value = 0.001972895328246894
$ws.Range("Q2").Value = 0.018673830595
$ws.Range("R2").Value = 0.168064475355
$ws.Range("S2").Value = 0.0002916492752488166
$ws.Range("T2").Value = 0.000302393504728266

# Row 3
$ws.Range("A3").Value = "ECs"
$ws.Range("B3").Value = "Il7"
$ws.Range("C3").Value = "Il7r"
$ws.Range("D3").Value = "Inflammatory-Mac"
$ws.Range("E3").Value = 1
$ws.Range("F3").Value = 0.3333333333333333
$ws.Range("G3").Value = 0.1621203333333333
$ws.Range("H3").Value = 0.486361
$ws.Range("I3").Value = 0.1481290732860657
$ws.Range("J3").Value = 0.1532739727236171
$ws.Range("K3").Value = 3
$ws.Range("L3").Value = 1
$ws.Range("M3").Value = 29.737678
$ws.Range("N3").Value = 89.21303400000001
$ws.Range("O3").Value = 0.5083135931924556
$ws.Range("P3").Value = 0.5093486651830572
$ws.Range("Q3").Value = 4.821082269919334
$ws.Range("R3").Value = 43.389740429274
$ws.Range("S3").Value = 0.07529602149830865
$ws.Range("T3").Value = 0.07806989341407868

# Row 4
$ws.Range("A4").Value = "ECs"
$ws.Range("B4").Value = "Il7"
$ws.Range("C4").Value = "Il7r"
$ws.Range("D4").Value = "MuSCs"
$ws.Range("E4").Value = 1
$ws.Range("F4").Value = 0.3333333333333333
$ws.Range("G4").Value = 0.1621203333333333
$ws.Range("H4").Value = 0.486361
$ws.Range("I4").Value = 0.1481290732860657
$ws.Range("J4").Value = 0.1532739727236171
$ws.Range("K4").Value = 1
$ws.Range("L4").Value = 0.5
$ws.Range("M4").Value = 0.356658
$ws.Range("N4").Value = 0.713316
$ws.Range("O4").Value = 0.00609644470293998
$ws.Range("P4").Value = 0.004072572539722363
$ws.Range("Q4").Value = 0.057821513846
$ws.Range("R4").Value = 0.346929083076
$ws.Range("S4").Value = 0.0009030607041862433
$ws.Range("T4").Value = 0.0006242193723683574

# Row 5
$ws.Range("A5").Value = "ECs"
$ws.Range("B5").Value = "Il7"
$ws.Range("C5").Value = "Il7r"
$ws.Range("D5").Value = "Resolving-Mac"
$ws.Range("E5").Value = 1
$ws.Range("F5").Value = 0.3333333333333333
$ws.Range("G5").Value = 0.1621203333333333
$ws.Range("H5").Value = 0.486361
$ws.Range("I5").Value = 0.1481290732860657
$ws.Range("J5").Value = 0.1532739727236171
$ws.Range("K5").Value = 3
$ws.Range("L5").Value = 1
$ws.Range("M5").Value = 28.293101
$ws.Range("N5").Value = 84.87930299999999
$ws.Range("O5").Value = 0.4836210759921153
$ws.Range("P5").Value = 0.4846058669489736
$ws.Range("Q5").Value = 4.586886965153666
$ws.Range("R5").Value = 41.281982686383
$ws.Range("S5").Value = 0.071638341808322
$ws.Range("T5").Value = 0.07427746643244178

# Row 6
$ws.Range("A6").Value = "FAPs"
$ws.Range("B6").Value = "Il7"
$ws.Range("C6").Value = "Il7r"
$ws.Range("D6").Value = "ECs"
$ws.Range("E6").Value = 3
$ws.Range("F6").Value = 1
$ws.Range("G6").Value = 0.8221213333333334
$ws.Range("H6").Value = 2.466364
$ws.Range("I6").Value = 0.7511708663032484
$ws.Range("J6").Value = 0.777260940870076
$ws.Range("K6").Value = 1
$ws.Range("L6").Value = 0.3333333333333333
$ws.Range("M6").Value = 0.115185
$ws.Range("N6").Value = 0.345555
$ws.Range("O6").Value = 0.001968886112489112
$ws.Range("P6").Value = 0.001972895328246894
$ws.Range("Q6").Value = 0.09469604578
$ws.Range("R6").Value = 0.85226441202
$ws.Range("S6").Value = 0.001478969886770881
$ws.Range("T6").Value = 0.001533454479071358

# Row 7
$ws.Range("A7").Value = "FAPs"
$ws.Range("B7").Value = "Il7"
$ws.Range("C7").Value = "Il7r"
$ws.Range("D7").Value = "Inflammatory-Mac"
$ws.Range("E7").Value = 3
$ws.Range("F7").Value = 1
$ws.Range("G7").Value = 0.8221213333333334
$ws.Range("H7").Value = 2.466364
$ws.Range("I7").Value = 0.7511708663032484
$ws.Range("J7").Value = 0.777260940870076
$ws.Range("K7").Value = 3
$ws.Range("L7").Value = 1
$ws.Range("M7").Value = 29.737678
$ws.Range("N7").Value = 89.21303400000001
$ws.Range("O7").Value = 0.5083135931924556
$ws.Range("P7").Value = 0.5093486651830572
$ws.Range("Q7").Value = 24.44797948759734
$ws.Range("R7").Value = 220.031815388376
$ws.Range("S7").Value = 0.3818303621520939
$ws.Range("T7").Value = 0.3958968227311004

# Row 8
$ws.Range("A8").Value = "FAPs"
$ws.Range("B8").Value = "Il7"
$ws.Range("C8").Value = "Il7r"
$ws.Range("D8").Value = "MuSCs"
$ws.Range("E8").Value = 3
$ws.Range("F8").Value = 1
$ws.Range("G8").Value = 0.8221213333333334
$ws.Range("H8").Value = 2.466364
$ws.Range("I8").Value = 0.7511708663032484
$ws.Range("J8").Value = 0.777260940870076
$ws.Range("K8").Value = 1
$ws.Range("L8").Value = 0.5
$ws.Range("M8").Value = 0.356658
$ws.Range("N8").Value = 0.713316
$ws.Range("O8").Value = 0.00609644470293998
$ws.Range("P8").Value = 0.004072572539722363
$ws.Range("Q8").Value = 0.293216150504
$ws.Range("R8").Value = 1.759296903024
$ws.Range("S8").Value = 0.004579471648877275
$ws.Range("T8").Value = 0.003165451563986239

# Row 9
$ws.Range("A9").Value = "FAPs"
$ws.Range("B9").Value = "Il7"
$ws.Range("C9").Value = "Il7r"
$ws.Range("D9").Value = "Resolving-Mac"
$ws.Range("E9").Value = 3
$ws.Range("F9").Value = 1
$ws.Range("G9").Value = 0.8221213333333334
$ws.Range("H9").Value = 2.466364
$ws.Range("I9").Value = 0.7511708663032484
$ws.Range("J9").Value = 0.777260940870076
$ws.Range("K9").Value = 3
$ws.Range("L9").Value = 1
$ws.Range("M9").Value = 28.293101
$ws.Range("N9").Value = 84.87930299999999
$ws.Range("O9").Value = 0.4836210759921153
$ws.Range("P9").Value = 0.4846058669489736
$ws.Range("Q9").Value = 23.26036191825467
$ws.Range("R9").Value = 209.343257264292
$ws.Range("S9").Value = 0.3632820626155063
$ws.Range("T9").Value = 0.3766652120959181

# Row 10
$ws.Range("A10").Value = "MuSCs"
$ws.Range("B10").Value = "Il7"
$ws.Range("C10").Value = "Il7r"
$ws.Range("D10").Value = "ECs"
$ws.Range("E10").Value = 1
$ws.Range("F10").Value = 0.5
$ws.Range("G10").Value = 0.1102115
$ws.Range("H10").Value = 0.220423
$ws.Range("I10").Value = 0.1007000604106861
$ws.Range("J10").Value = 0.06946508640630693
$ws.Range("K10").Value = 1
$ws.Range("L10").Value = 0.3333333333333333
$ws.Range("M10").Value = 0.115185
$ws.Range("N10").Value = 0.345555
$ws.Range("O10").Value = 0.001968886112489112
$ws.Range("P10").Value = 0.001972895328246894
$ws.Range("Q10").Value = 0.0126947116275
$ws.Range("R10").Value = 0.076168269765
$ws.Range("S10").Value = 0.0001982669504694144
$ws.Range("T10").Value = 0.0001370473444472698

# Row 11
$ws.Range("A11").Value = "MuSCs"
$ws.Range("B11").Value = "Il7"
$ws.Range("C11").Value = "Il7r"
$ws.Range("D11").Value = "Inflammatory-Mac"
$ws.Range("E11").Value = 1
$ws.Range("F11").Value = 0.5
$ws.Range("G11").Value = 0.1102115
$ws.Range("H11").Value = 0.220423
$ws.Range("I11").Value = 0.1007000604106861
$ws.Range("J11").Value = 0.06946508640630693
$ws.Range("K11").Value = 3
$ws.Range("L11").Value = 1
$ws.Range("M11").Value = 29.737678
$ws.Range("N11").Value = 89.21303400000001
$ws.Range("O11").Value = 0.5083135931924556
$ws.Range("P11").Value = 0.5093486651830572
$ws.Range("Q11").Value = 3.277434098897
$ws.Range("R11").Value = 19.664604593382
$ws.Range("S11").Value = 0.05118720954205318
$ws.Range("T11").Value = 0.03538194903787817

# Row 12
$ws.Range("A12").Value = "MuSCs"
$ws.Range("B12").Value = "Il7"
$ws.Range("C12").Value = "Il7r"
$ws.Range("D12").Value = "MuSCs"
$ws.Range("E12").Value = 1
$ws.Range("F12").Value = 0.5
$ws.Range("G12").Value = 0.1102115
$ws.Range("H12").Value = 0.220423
$ws.Range("I12").Value = 0.1007000604106861
$ws.Range("J12").Value = 0.06946508640630693
$ws.Range("K12").Value = 1
$ws.Range("L12").Value = 0.5
$ws.Range("M12").Value = 0.356658
$ws.Range("N12").Value = 0.713316
$ws.Range("O12").Value = 0.00609644470293998
$ws.Range("P12").Value = 0.004072572539722363
$ws.Range("Q12").Value = 0.039307813167
$ws.Range("R12").Value = 0.157231252668
$ws.Range("S12").Value = 0.0006139123498764631
$ws.Range("T12").Value = 0.0002829016033677668

# Row 13
$ws.Range("A13").Value = "MuSCs"
$ws.Range("B13").Value = "Il7"
$ws.Range("C13").Value = "Il7r"
$ws.Range("D13").Value = "Resolving-Mac"
$ws.Range("E13").Value = 1
$ws.Range("F13").Value = 0.5
$ws.Range("G13").Value = 0.1102115
$ws.Range("H13").Value = 0.220423
$ws.Range("I13").Value = 0.1007000604106861
$ws.Range("J13").Value = 0.06946508640630693
$ws.Range("K13").Value = 3
$ws.Range("L13").Value = 1
$ws.Range("M13").Value = 28.293101
$ws.Range("N13").Value = 84.87930299999999
$ws.Range("O13").Value = 0.4836210759921153
$ws.Range("P13").Value = 0.4846058669489736
$ws.Range("Q13").Value = 3.1182251008615
$ws.Range("R13").Value = 18.709350605169
$ws.Range("S13").Value = 0.04870067156828701
$ws.Range("T13").Value = 0.03366318842061373
